$d = $word.ActiveDocument
$t = $d.Tables(1)

$t.Cell(1, 1).Range.Text = "Rank"
$t.Cell(1, 2).Range.Text = "Description"
$t.Cell(1, 3).Range.Text = "Last Week"
$t.Cell(1, 4).Range.Text = "Weeks on List"
